# Update odds values on Sheet1 (rows 2, 4, 5) to match the latest
# FlashScore snapshot for 2025-03-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Union de Santa Fe vs Racing Club) ---
$ws.Range("I2").Value  = 2.4
$ws.Range("J2").Value  = 3.75
$ws.Range("M2").Value  = 1.1
$ws.Range("N2").Value  = 7
$ws.Range("AA2").Value = 1.95
$ws.Range("AB2").Value = 1.8
$ws.Range("AC2").Value = 8
$ws.Range("AO2").Value = 11

# --- Row 4 (Montevideo City vs Juventud) ---
$ws.Range("G4").Value  = 1.91
$ws.Range("I4").Value  = 3.6
$ws.Range("J4").Value  = 2.75
$ws.Range("L4").Value  = 4.5
$ws.Range("S4").Value  = 2.2
$ws.Range("T4").Value  = 1.65
$ws.Range("AJ4").Value = 7
$ws.Range("AN4").Value = 9
$ws.Range("AO4").Value = 17

# --- Row 5 (Penarol vs Liverpool M.) ---
$ws.Range("G5").Value  = 1.57
$ws.Range("I5").Value  = 5.5
$ws.Range("J5").Value  = 2.25
$ws.Range("M5").Value  = 1.08
$ws.Range("N5").Value  = 8
$ws.Range("S5").Value  = 2.15
$ws.Range("T5").Value  = 1.67
$ws.Range("W5").Value  = 4
$ws.Range("X5").Value  = 1.22
$ws.Range("AC5").Value = 5.5
$ws.Range("AD5").Value = 6.5
$ws.Range("AE5").Value = 9
$ws.Range("AF5").Value = 11
$ws.Range("AI5").Value = 8
$ws.Range("AK5").Value = 21
$ws.Range("AO5").Value = 29
$ws.Range("AP5").Value = 19
$ws.Range("AQ5").Value = 67
$ws.Range("AR5").Value = 51
